$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9817726016044617
$ws.Range("B1").Value = 1.520613074302673
$ws.Range("C1").Value = 6.654296875
$ws.Range("D1").Value = 2.122264623641968
$ws.Range("E1").Value = 0.9412879347801208
